$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct classification value typo: 2023 (not in AAA-Objektartenkatalog) -> 2013
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = 2013

# Reflect the selection state saved in the workbook
[void]$ws.Range("A13").Select()
